# fix import authentication error!
# Rename the Chinese header labels in row 1 to match the new (generic
# data-type) column headers, and fix the typo'd key names in the sample
# "object" cell G3 (datatype/datavalue -> mode/value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "数字"
$ws.Range("B1").Value = "字符串"
$ws.Range("C1").Value = "布尔"
$ws.Range("D1").Value = "数字数组"
$ws.Range("E1").Value = "字符串数组"
$ws.Range("F1").Value = "布尔数组"
$ws.Range("G1").Value = "对象"
# H1 ("对象数组") is unchanged.

$ws.Range("G3").Value = "mode:percent;value:45"

# The user resized several columns (A, B, D, F, G) by hand after retyping
# the shorter headers; C, E and H were left at their existing best-fit size.
$ws.Columns.Item(1).ColumnWidth = 4.2857142857142856
$ws.Columns.Item(2).ColumnWidth = 11.160714285714286
$ws.Columns.Item(4).ColumnWidth = 9.9107142857142856
$ws.Columns.Item(6).ColumnWidth = 12.285714285714286
$ws.Columns.Item(7).ColumnWidth = 21.410714285714286

# Final cursor position left on D11.
$ws.Range("D11").Select()
